{"js": "// Fix the \"missing times\" issue: update cumulative_time_error / gap_error\n// values in the Spain GP comparison table. Each old numeric value is\n// unique in the document, so we can safely search-and-replace each one\n// with its corrected counterpart.\nconst replacements = [\n  [\"204.792\", \"208.570\"],\n  [\"-2.219\", \"-1.739\"],\n  [\"207.885\", \"211.184\"],\n  [\"197.466\", \"201.480\"],\n  [\"-9.545\", \"-8.829\"],\n  [\"203.943\", \"206.707\"],\n  [\"-3.068\", \"-3.602\"],\n  [\"201.273\", \"205.863\"],\n  [\"-5.738\", \"-4.446\"],\n  [\"210.740\", \"214.977\"],\n  [\"3.729\", \"4.668\"],\n  [\"202.585\", \"207.509\"],\n  [\"-4.426\", \"-2.800\"],\n  [\"189.129\", \"193.801\"],\n  [\"-17.882\", \"-16.509\"],\n  [\"192.090\", \"197.418\"],\n  [\"-14.921\", \"-12.892\"],\n  [\"204.076\", \"209.184\"],\n  [\"-2.935\", \"-1.125\"],\n  [\"211.002\", \"216.674\"],\n  [\"3.991\", \"6.365\"],\n  [\"283.465\", \"289.006\"],\n  [\"-3.320\", \"-1.078\"],\n  [\"279.443\", \"285.251\"],\n  [\"-7.342\", \"-4.832\"],\n  [\"253.380\", \"259.042\"],\n  [\"-33.405\", \"-31.042\"],\n  [\"276.298\", \"283.002\"],\n  [\"-10.488\", \"-7.081\"],\n  [\"308.946\", \"315.211\"],\n  [\"22.161\", \"25.127\"],\n  [\"270.860\", \"277.194\"],\n  [\"-15.925\", \"-12.889\"],\n  [\"268.583\", \"279.218\"],\n  [\"-18.202\", \"-10.866\"],\n  [\"285.361\", \"291.305\"],\n  [\"-1.424\", \"1.222\"],\n  [\"356.036\", \"362.479\"],\n  [\"-10.216\", \"-7.071\"],\n];\n\nfor (const [oldValue, newValue] of replacements) {\n  const results = context.document.body.search(oldValue, {\n    matchCase: true,\n    matchWholeWord: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newValue, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Fix the \"missing times\" issue: update cumulative_time_error / gap_error\n# values in the Spain GP comparison table. Each old numeric value is\n# unique in the document, so Find/Replace on the whole document content\n# can safely target each one individually.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"204.792\", \"208.570\"),\n    @(\"-2.219\", \"-1.739\"),\n    @(\"207.885\", \"211.184\"),\n    @(\"197.466\", \"201.480\"),\n    @(\"-9.545\", \"-8.829\"),\n    @(\"203.943\", \"206.707\"),\n    @(\"-3.068\", \"-3.602\"),\n    @(\"201.273\", \"205.863\"),\n    @(\"-5.738\", \"-4.446\"),\n    @(\"210.740\", \"214.977\"),\n    @(\"3.729\", \"4.668\"),\n    @(\"202.585\", \"207.509\"),\n    @(\"-4.426\", \"-2.800\"),\n    @(\"189.129\", \"193.801\"),\n    @(\"-17.882\", \"-16.509\"),\n    @(\"192.090\", \"197.418\"),\n    @(\"-14.921\", \"-12.892\"),\n    @(\"204.076\", \"209.184\"),\n    @(\"-2.935\", \"-1.125\"),\n    @(\"211.002\", \"216.674\"),\n    @(\"3.991\", \"6.365\"),\n    @(\"283.465\", \"289.006\"),\n    @(\"-3.320\", \"-1.078\"),\n    @(\"279.443\", \"285.251\"),\n    @(\"-7.342\", \"-4.832\"),\n    @(\"253.380\", \"259.042\"),\n    @(\"-33.405\", \"-31.042\"),\n    @(\"276.298\", \"283.002\"),\n    @(\"-10.488\", \"-7.081\"),\n    @(\"308.946\", \"315.211\"),\n    @(\"22.161\", \"25.127\"),\n    @(\"270.860\", \"277.194\"),\n    @(\"-15.925\", \"-12.889\"),\n    @(\"268.583\", \"279.218\"),\n    @(\"-18.202\", \"-10.866\"),\n    @(\"285.361\", \"291.305\"),\n    @(\"-1.424\", \"1.222\"),\n    @(\"356.036\", \"362.479\"),\n    @(\"-10.216\", \"-7.071\")\n)\n\nforeach ($pair in $replacements) {\n    $oldValue = $pair[0]\n    $newValue = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldValue\n    $find.Replacement.Text = $newValue\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n}\n"}
